$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original layout: A=idx, B=rec_yds, C=rec_td, D=fumbles, E=fantasy points (17 data rows, rows 2-17)
# New layout:       A=idx, B=rec_yds, C=rec_td, D=fumbles, E=height, F=weight, G=fantasy points

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row() + $usedRange.Rows.Count() - 1

# Preserve the existing "fantasy points" header text and column values before overwriting column E
$fpHeader = $ws.Range("E1").Value()
$fpValues = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $fpValues[$r] = $ws.Cells.Item($r, 5).Value()
}

# Give the two new header cells (F1, G1) the same look as the existing bordered/bold header cells
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)

# New header row: E becomes "height", F becomes "weight", G gets the old "fantasy points" header text
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = $fpHeader

# Fill new height (E) / weight (F) constant values, and shift the old fantasy-points values to G
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.5
    $ws.Cells.Item($r, 6).Value = 260
    $ws.Cells.Item($r, 7).Value = $fpValues[$r]
}
